# Generate Report for Handback
#
# Refresh the "Correspond Handoff Datetime" (col E) and "Correspond
# Handback DateTime" (col H) timestamps for the handback-status report.
# Rows 3 and 4 of each language sheet shared the same timestamp text in
# the original report, so both rows move together (matches the source
# xlsx, where the two rows pointed at the same shared-string entry).

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-21 00:19:36"
$wsZh.Range("E4").Value = "2016-03-21 00:19:36"
$wsZh.Range("H3").Value = "2016-03-21 00:19:58"
$wsZh.Range("H4").Value = "2016-03-21 00:19:58"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-21 00:19:39"
$wsDe.Range("E4").Value = "2016-03-21 00:19:39"
$wsDe.Range("H3").Value = "2016-03-21 00:20:05"
$wsDe.Range("H4").Value = "2016-03-21 00:20:05"
